# The edit shuffles the data rows (rows 2-33) of the single worksheet into a
# new order. Row 1 (header) and row 30 stay where they are; every other row's
# entire content (columns A:AY) is relocated to a different row per the
# mapping below (new row number -> old row number that supplies its data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 33
$lastCol = "AY"

# new row -> source (old) row
$map = @{
    2  = 7
    3  = 14
    4  = 21
    5  = 8
    6  = 18
    7  = 13
    8  = 3
    9  = 12
    10 = 16
    11 = 2
    12 = 10
    13 = 4
    14 = 15
    15 = 32
    16 = 24
    17 = 9
    18 = 33
    19 = 28
    20 = 27
    21 = 23
    22 = 6
    23 = 19
    24 = 11
    25 = 5
    26 = 29
    27 = 26
    28 = 22
    29 = 25
    30 = 30
    31 = 20
    32 = 17
    33 = 31
}

# Columns that hold text (as opposed to real numbers/booleans) in this sheet.
# Several of them contain digit- or date-looking strings (e.g. "1",
# "2023-08-28", "00:00") that Excel would otherwise auto-convert to a number
# or date serial on assignment. Forcing the Text number format first makes
# the write preserve the original string type.
$textCols = @("C","D","F","G","H","I","J","K","L","M","N","O","P", `
              "T","U","V","W","X","Y","Z","AA","AB","AC", `
              "AF","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR", `
              "AS","AU","AW","AX","AY")
foreach ($col in $textCols) {
    $ws.Range(($col + $firstDataRow + ":" + $col + $lastDataRow)).NumberFormat = "@"
}

# 1) Snapshot every source row's full contents BEFORE any writes happen,
#    since rows are both sources and destinations (an in-place permutation).
$snapshot = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rng = $ws.Range("A" + $r + ":" + $lastCol + $r)
    $snapshot[$r] = $rng.Value2
}

# 2) Write each new row from its recorded source snapshot.
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $srcRow = $map[$r]
    $dst = $ws.Range("A" + $r + ":" + $lastCol + $r)
    $dst.Value = $snapshot[$srcRow]
}

Write-Output "row shuffle applied"
